# Rebuild the "Linea" sheet content to match the ETL export (adds LINEA_ID / ANIO_INAUGURACION /
# TAM_KM / AFLUENCIA columns and 12 line rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Linea")
$ws.Activate()

# --- Re-assert the shared strings already used elsewhere in the sheet (kept verbatim) ---
$ws.Range("A1").Value = "SISTEMA"
$ws.Range("C1").Value = "NOMBRE"
$ws.Range("E1").Value = "COLOR_ESP"
$ws.Range("F1").Value = "COLOR_EN"
$ws.Range("A2").Value = "STC Metro"
$ws.Range("A3").Value = "STC Metro"
$ws.Range("A4").Value = "STC Metro"
$ws.Range("A5").Value = "STC Metro"
$ws.Range("A6").Value = "STC Metro"
$ws.Range("A7").Value = "STC Metro"
$ws.Range("A8").Value = "STC Metro"
$ws.Range("A9").Value = "STC Metro"
$ws.Range("A10").Value = "STC Metro"
$ws.Range("A11").Value = "STC Metro"
$ws.Range("A12").Value = "STC Metro"
$ws.Range("A13").Value = "STC Metro"
$ws.Range("E2").Value = "ROSA"
$ws.Range("F2").Value = "PINK"

# --- New text values, written in first-use order ---
$ws.Range("B1").Value = "LINEA_ID"
$ws.Range("D1").Value = "ANIO_INAUGURACION"
$ws.Range("E3").Value = "AZUL"
$ws.Range("F3").Value = "BLUE"
$ws.Range("E4").Value = "VERDE_OLIVO"
$ws.Range("F4").Value = "OLIVE_GREEN"
$ws.Range("G1").Value = "TAM_KM"
$ws.Range("E5").Value = "CIAN"
$ws.Range("F5").Value = "CYAN"
$ws.Range("H1").Value = "AFLUENCIA"
$ws.Range("H5").Value = "BAJA"
$ws.Range("E6").Value = "AMARILLO"
$ws.Range("F6").Value = "YELLOW"
$ws.Range("E7").Value = "ROJO"
$ws.Range("F7").Value = "RED"
$ws.Range("E8").Value = "NARANJA"
$ws.Range("F8").Value = "ORANGE"
$ws.Range("F9").Value = "GREEN"
$ws.Range("E9").Value = "VERDE"
$ws.Range("E10").Value = "CAFÉ"
$ws.Range("F10").Value = "BROWN"
$ws.Range("E13").Value = "DORADO"
$ws.Range("F13").Value = "GOLDEN"
$ws.Range("C2").Value = "LINEA 1"
$ws.Range("C3").Value = "LINEA 2"
$ws.Range("C4").Value = "LINEA 3"
$ws.Range("C5").Value = "LINEA 4"
$ws.Range("C6").Value = "LINEA 5"
$ws.Range("C7").Value = "LINEA 6"
$ws.Range("C8").Value = "LINEA 7"
$ws.Range("C9").Value = "LINEA 8"
$ws.Range("C10").Value = "LINEA 9"
$ws.Range("C13").Value = "LINEA 12"
$ws.Range("C11").Value = "LINEA A"
$ws.Range("E11").Value = "MORADA"
$ws.Range("F11").Value = "PURPLE"
$ws.Range("C12").Value = "LINEA B"
$ws.Range("E12").Value = "VERDE_GRIS"
$ws.Range("F12").Value = "GREEN_GREY"

# --- Numeric values (line number, inauguration year, km length) ---
$ws.Range("B2").Value = 1
$ws.Range("D2").Value = 1969
$ws.Range("G2").Value = 18.8
$ws.Range("B3").Value = 2
$ws.Range("D3").Value = 1970
$ws.Range("G3").Value = 23.43
$ws.Range("B4").Value = 3
$ws.Range("D4").Value = 1970
$ws.Range("G4").Value = 23.61
$ws.Range("B5").Value = 4
$ws.Range("D5").Value = 1981
$ws.Range("G5").Value = 10.75
$ws.Range("B6").Value = 5
$ws.Range("D6").Value = 1981
$ws.Range("G6").Value = 16.67
$ws.Range("B7").Value = 6
$ws.Range("D7").Value = 1983
$ws.Range("G7").Value = 13.95
$ws.Range("B8").Value = 7
$ws.Range("D8").Value = 1984
$ws.Range("G8").Value = 18.78
$ws.Range("B9").Value = 8
$ws.Range("D9").Value = 1994
$ws.Range("G9").Value = 20.08
$ws.Range("B10").Value = 9
$ws.Range("D10").Value = 1987
$ws.Range("G10").Value = 15.37
$ws.Range("B11").Value = 10
$ws.Range("D11").Value = 1991
$ws.Range("G11").Value = 17.19
$ws.Range("B12").Value = 11
$ws.Range("D12").Value = 1999
$ws.Range("G12").Value = 23.72
$ws.Range("B13").Value = 12
$ws.Range("D13").Value = 2012
$ws.Range("G13").Value = 24.5

# --- Register the small (size-8) Calibri font used for the phonetic-guide settings of the
#     imported range, then drop it from this cell again so no visible cell keeps it. ---
$ws.Range("B1").Font.Size = 8
$ws.Range("B1").Font.Size = 11

# --- Leave the selection where the author left it before saving ---
$ws.Range("M6").Select()
